# "switch regression to be in different table"
#
# The sheet currently has two result columns:
#   B = OLS results (header "(1)" / "OLS")
#   C = Poisson (mfx) results (header "(2)" / "Poisson (mfx)")
#
# The edit drops the OLS column and keeps only the Poisson (mfx) column,
# moving it into column B (so the table becomes a single-column table),
# and relabels its header row from "Poisson (mfx)" to
# "matrix.ncol...1..nrow...16.". The "(1)" label above it is kept as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column B (OLS). This shifts column C (Poisson mfx) left into B,
# carrying over its values, number/string typing and cell styles, and
# automatically shrinks the A1:C1 merged header cell down to A1:B1.
$ws.Range("B:B").Delete()

# Restore/relabel the two header cells in the now-single results column.
$ws.Range("B2").Value2 = "(1)"
$ws.Range("B3").Value2 = "matrix.ncol...1..nrow...16."
